$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.1125
$ws.Range("C6").Value = -12.864
$ws.Range("C7").Value = -12.46569999999999
$ws.Range("D7").Value = -7.474499999999995
$ws.Range("C8").Value = -13.3273
$ws.Range("D11").Value = -7.744799999999997
$ws.Range("D12").Value = -7.260899999999999
$ws.Range("D15").Value = -8.88619999999999
$ws.Range("C16").Value = -13.93939999999999
$ws.Range("C20").Value = -12.125
$ws.Range("D20").Value = -7.987599999999994
$ws.Range("C21").Value = -12.4039
$ws.Range("D21").Value = -8.006099999999998
$ws.Range("D22").Value = -8.132500000000006
$ws.Range("D23").Value = -7.415899999999996
$ws.Range("C28").Value = -12.0509
$ws.Range("C29").Value = -11.07780000000001
$ws.Range("D29").Value = -7.152099999999997
$ws.Range("C30").Value = -12.6863
$ws.Range("C32").Value = -12.7874
$ws.Range("D34").Value = -7.754100000000001
$ws.Range("C40").Value = -11.93960000000001
$ws.Range("D42").Value = -8.541400000000005
$ws.Range("D43").Value = -8.202700000000002
$ws.Range("D44").Value = -7.023300000000002
$ws.Range("D45").Value = -7.393599999999997
$ws.Range("C46").Value = -13.80099999999999
$ws.Range("D46").Value = -8.5991
$ws.Range("D50").Value = -8.013499999999995
$ws.Range("C51").Value = -12.4683
$ws.Range("D51").Value = -7.846699999999996
$ws.Range("C52").Value = -11.1008
$ws.Range("C57").Value = -13.71219999999999
$ws.Range("D57").Value = -8.713000000000001
$ws.Range("C59").Value = -12.65229999999999
$ws.Range("C62").Value = -13.63239999999998
$ws.Range("D65").Value = -7.849399999999998
$ws.Range("C66").Value = -11.09140000000001
$ws.Range("D66").Value = -7.295899999999998
$ws.Range("D67").Value = -6.380799999999999
$ws.Range("C73").Value = -11.22930000000001
$ws.Range("C74").Value = -12.10900000000001
$ws.Range("C77").Value = -12.2357
$ws.Range("D79").Value = -6.262100000000003
$ws.Range("D84").Value = -8.823800000000002
$ws.Range("D87").Value = -7.892599999999995
$ws.Range("C92").Value = -11.3886
$ws.Range("D92").Value = -6.624500000000004
$ws.Range("D97").Value = -8.490500000000003
$ws.Range("C100").Value = -12.1569
